$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcCols = @("B","C","D","E","F","G","H","I","J","K")
$dstCols = @("A","B","C","D","E","F","G","H","I","J")

for ($r = 4; $r -le 7; $r++) {
    # Read all source values for the row first so overwriting earlier
    # columns does not clobber values still to be read.
    $rowvals = @()
    foreach ($sc in $srcCols) {
        $rowvals += ,$ws.Range("$sc$r").Value2
    }

    for ($i = 0; $i -lt $dstCols.Length; $i++) {
        $dc = $dstCols[$i]
        $v = $rowvals[$i]
        $destCell = $ws.Range("$dc$r")

        if ($v -eq "TRUE" -or $v -eq "FALSE") {
            # A bare assignment of the literal text TRUE/FALSE gets
            # auto-coerced to a native boolean cell. Forcing a text
            # entry with a leading apostrophe keeps it a normal string,
            # then we restore the row's usual (non quote-prefixed) look
            # by re-pasting the formatting from column A of the same row.
            $destCell.Value = "'" + $v
            $ws.Range("A$r").Copy()
            $destCell.PasteSpecial(-4122)
        } else {
            $destCell.Value = $v
        }
    }
}

$excel.CutCopyMode = 0

# The former last column (K) is now unused for these rows - remove it
# completely instead of leaving blank/styled cells behind.
$ws.Range("K4:K7").Clear()

$ws.Range("A4").Select()
